$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing text content first (this keeps every cell's style, each
# row's custom height, and the column widths intact -- ClearContents only
# touches values) so the shared-string table rebuilds from empty and is
# repopulated strictly in the order the cells are (re)written below, the
# same left-to-right/top-to-bottom order the source file itself uses.
$ws.Range("A1:C22").ClearContents()

$ws.Range("B1").Value = 'Ementa atual:'

$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'

$ws.Range("B2").Value = 'LOQ4233'

$ws.Range("C2").Value = 'LOQ4233'

$ws.Range("A3").Value = 'Nome:'

$ws.Range("B3").Value = ' Gestão de Negócios'

$ws.Range("C3").Value = ' Gestão de Negócios'

$ws.Range("A4").Value = 'Name:'

$ws.Range("B4").Value = 'Business Management'

$ws.Range("C4").Value = 'Business Management'

$ws.Range("A5").Value = 'Créditos-aula:'

$ws.Range("B5").Formula = '="2"'
$ws.Range("B5").Copy()
$ws.Range("B5").PasteSpecial(-4163)

$ws.Range("C5").Formula = '="2"'
$ws.Range("C5").Copy()
$ws.Range("C5").PasteSpecial(-4163)

$ws.Range("A6").Value = 'Créditos-trabalho'

$ws.Range("B6").Formula = '="0"'
$ws.Range("B6").Copy()
$ws.Range("B6").PasteSpecial(-4163)

$ws.Range("C6").Formula = '="0"'
$ws.Range("C6").Copy()
$ws.Range("C6").PasteSpecial(-4163)

$ws.Range("A7").Value = 'Carga horária:'

$ws.Range("B7").Value = '30 h'

$ws.Range("C7").Value = '30 h'

$ws.Range("A8").Value = 'Ativação:'

$ws.Range("B8").Formula = '="01/01/2022"'
$ws.Range("B8").Copy()
$ws.Range("B8").PasteSpecial(-4163)

$ws.Range("C8").Formula = '="01/01/2022"'
$ws.Range("C8").Copy()
$ws.Range("C8").PasteSpecial(-4163)

$ws.Range("A9").Value = 'Semestre ideal:'

$ws.Range("B9").Value = 'EF-7,EM-6,EA-5,EB-4,EQD-7,EQN-11'

$ws.Range("C9").Value = 'EF-7,EM-6,EA-5,EB-4,EQD-7,EQN-11'

$ws.Range("A10").Value = 'Objetivos:'

$ws.Range("B10").Value = 'Apresentar ao aluno o conceito de uma organização e os fundamentos de sua administração;Caracterizar as diversas áreas funcionais existentes nas organizações;Despertar o interesse dos alunos para questões de gestão'

$ws.Range("C10").Value = 'Apresentar ao aluno o conceito de uma organização e os fundamentos de sua administração;Caracterizar as diversas áreas funcionais existentes nas organizações;Despertar o interesse dos alunos para questões de gestão'

$ws.Range("A11").Value = 'Objectives:'

$ws.Range("B11").Value = 'To present to the student the concept of an organization and the foundations of its administration; to characterize the various functional areas existing in the organizations; to awaken the interest of the students for management issues.'
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)

$ws.Range("C11").Value = 'To present to the student the concept of an organization and the foundations of its administration; to characterize the various functional areas existing in the organizations; to awaken the interest of the students for management issues.'
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("A12").Value = 'Docentes responsáveis:'

$ws.Range("B13").Value = '849935 - Humberto Felipe da Silva'

$ws.Range("C13").Value = '849935 - Humberto Felipe da Silva'

$ws.Range("A14").Value = 'Programa resumido:'

$ws.Range("B14").Value = '1 - A Administração das Organizações. 2 - O processo administrativo. 3 – Processos de Gestão'

$ws.Range("C14").Value = '1 - A Administração das Organizações. 2 - O processo administrativo. 3 – Processos de Gestão'

$ws.Range("A15").Value = 'Short syllabus:'

$ws.Range("B15").Value = 'The Administration of Organizations. 2 - The Administrative Process. 3 - Management Processes'
$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("C15").Value = 'The Administration of Organizations. 2 - The Administrative Process. 3 - Management Processes'
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("A16").Value = 'Programa:'

$ws.Range("B16").Value = '1 - A Administração das organizações - definindo a administração2 - O processo administrativo: planejamento, organização, direção, controle3 – Processos de Gestão: Marketing, Finanças, Gestão de Pessoas, Produção e Operações, Pesquisa e Desenvolvimento, Tecnologia da Informação, Logística e Meio Ambiente.'

$ws.Range("C16").Value = '1 - A Administração das organizações - definindo a administração2 - O processo administrativo: planejamento, organização, direção, controle3 – Processos de Gestão: Marketing, Finanças, Gestão de Pessoas, Produção e Operações, Pesquisa e Desenvolvimento, Tecnologia da Informação, Logística e Meio Ambiente.'

$ws.Range("A17").Value = 'Syllabus:'

$ws.Range("B17").Value = '- The Administration of organizations - defining the administration 2 - The administrative process: planning, organization, direction, control 3 - Management Processes: Marketing, Finance, People Management, Production and Operations, Research and Development, Information Technology, Logistics and Environment.'
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)

$ws.Range("C17").Value = '- The Administration of organizations - defining the administration 2 - The administrative process: planning, organization, direction, control 3 - Management Processes: Marketing, Finance, People Management, Production and Operations, Research and Development, Information Technology, Logistics and Environment.'
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)

$ws.Range("A18").Value = 'Avaliação:'

$ws.Range("A19").Value = 'Método:'

$ws.Range("B19").Value = 'O sistema de avaliação será continuo com aplicação de avaliações escritas utilizando-se avaliações em papel como em sistema informacionais, ademais serão realizados seminários, projetos, entrega de trabalho em formato de artigo e Estudos de Casos.'

$ws.Range("C19").Value = 'O sistema de avaliação será continuo com aplicação de avaliações escritas utilizando-se avaliações em papel como em sistema informacionais, ademais serão realizados seminários, projetos, entrega de trabalho em formato de artigo e Estudos de Casos.'

$ws.Range("A20").Value = 'Critério:'

$ws.Range("B20").Value = 'Avaliações em diversos formatos realizadas no decorrer do semestre. O peso maior da avaliação será aplicado ao Seminário Final da Disciplina, quando serão realizadas a apresentação oral do trabalho bem como a entrega do trabalho em formato de artigo; essa avaliação representará 70% da média do semestre.'

$ws.Range("C20").Value = 'Avaliações em diversos formatos realizadas no decorrer do semestre. O peso maior da avaliação será aplicado ao Seminário Final da Disciplina, quando serão realizadas a apresentação oral do trabalho bem como a entrega do trabalho em formato de artigo; essa avaliação representará 70% da média do semestre.'

$ws.Range("A21").Value = 'Norma de recuperação:'

$ws.Range("B21").Value = 'NF = (MF + PR)/ 2 , onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota do trabalho de recuperação'

$ws.Range("C21").Value = 'NF = (MF + PR)/ 2 , onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota do trabalho de recuperação'

$ws.Range("A22").Value = 'Bibliografia:'

$ws.Range("B22").Value = 'LEMOS, Paulo de Mattos et al. Gestão estratégica de empresas. Rio de Janeiro: Fundação Getúlio Vargas, 2014.Ludovico, Nelson. Gestão estratégica de negócios. São Paulo: Saraiva, 2018Serra, Fernando Ribeiro et al. Gestão estratégica: conceitos e casos. São Paulo: Atlas, 2014.'

$ws.Range("C22").Value = 'LEMOS, Paulo de Mattos et al. Gestão estratégica de empresas. Rio de Janeiro: Fundação Getúlio Vargas, 2014.Ludovico, Nelson. Gestão estratégica de negócios. São Paulo: Saraiva, 2018Serra, Fernando Ribeiro et al. Gestão estratégica: conceitos e casos. São Paulo: Atlas, 2014.'

$excel.CutCopyMode = $false